$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet6")

# Insert 11 blank rows before the existing data row (row 2), pushing it down to row 13.
$ws.Range("A2:A12").EntireRow.Insert()

# --- Row 2 (new) ---
$ws.Range("A2").Value = "56ac9c972233f1dc2aa17ec0"
$ws.Range("D2").Value = "---"
$ws.Range("E2").Value = "yes"
$ws.Range("L2").NumberFormat = "#,##0"
$ws.Range("V2").NumberFormat = "#,##0"

# --- Row 3 (new) ---
$ws.Range("A3").Value = "56ac9cdfa6e62d0029ab67b6"
$ws.Range("D3").Value = "---"
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 14766666666666600
$ws.Range("L3").NumberFormat = "#,##0"
$ws.Range("M3").Value = " 0 0 0"
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 3
$ws.Range("S3").Value = 3

# --- Row 4 (new) ---
$ws.Range("A4").Value = "56accf0c6cbabdfc1472b0dc"
$ws.Range("D4").Value = "---"
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 14766666666666600
$ws.Range("L4").NumberFormat = "#,##0"
$ws.Range("M4").Value = " 0 0 0"

# --- Row 5 (new) ---
$ws.Range("A5").Value = "56acd68759f671c40e4a9074"
$ws.Range("D5").Value = "---"
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 14766666666666600
$ws.Range("L5").NumberFormat = "#,##0"
$ws.Range("M5").Value = " 0 0 0"

# --- Row 6 (new) ---
$ws.Range("A6").Value = "56acd9938584f07025b14a71"
$ws.Range("D6").Value = "---"
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 14766666666666600
$ws.Range("L6").NumberFormat = "#,##0"
$ws.Range("M6").Value = " 0 0 0"

# --- Row 7 (new) ---
$ws.Range("A7").Value = "56acdbf659e7f7341cbf50fd"
$ws.Range("D7").Value = "---"
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 14766666666666600
$ws.Range("L7").NumberFormat = "#,##0"
$ws.Range("M7").Value = " 0 0 0"

# --- Row 8 (new) ---
$ws.Range("A8").Value = "56acdcbda3b8ebec07ae9daa"
$ws.Range("D8").Value = "---"
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 14766666666666600
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("M8").Value = " 0 0 0"

# --- Row 9 (new) ---
$ws.Range("A9").Value = "56acdefc11e52d64287f6d68"
$ws.Range("D9").Value = "---"
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 14766666666666600
$ws.Range("L9").NumberFormat = "#,##0"
$ws.Range("M9").Value = " 0 0 0"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "56ace0ccc7dc11f012e3da04"
$ws.Range("D10").Value = "---"
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 14766666666666600
$ws.Range("L10").NumberFormat = "#,##0"
$ws.Range("M10").Value = " 0 0 0"
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 2
$ws.Range("S10").Value = 2

# --- Row 11 (new) ---
$ws.Range("A11").Value = "56ace2005329f068279cfac3"
$ws.Range("D11").Value = "---"
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 14766666666666600
$ws.Range("L11").NumberFormat = "#,##0"
$ws.Range("M11").Value = " 0 0 0"

# --- Row 12 (new) ---
$ws.Range("A12").Value = "56ace3a08fe82410206d7e23"
$ws.Range("D12").Value = "---"
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 14766666666666600
$ws.Range("L12").NumberFormat = "#,##0"
$ws.Range("M12").Value = " 0 0 0"
$ws.Range("P12").Value = 2
$ws.Range("Q12").Value = 2
$ws.Range("S12").Value = 2

# --- Row 13 already holds the original row-2 data (shifted by the insert above) ---

# --- Row 14 (new) ---
$ws.Range("A14").Value = "56adf71bd9b0d1ac0d403bff"
$ws.Range("D14").Value = "---"
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 14766666666666600
$ws.Range("L14").NumberFormat = "#,##0"
$ws.Range("M14").Value = " 0 0 0"
$ws.Range("N14").Value = 282
$ws.Range("O14").Value = 6
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 2
$ws.Range("S14").Value = 2
$ws.Range("V14").Value = 20258289768529800
$ws.Range("V14").NumberFormat = "#,##0"
$ws.Range("W14").Value = 21
$ws.Range("X14").Value = 2
$ws.Range("Y14").Value = 2
$ws.Range("AA14").Value = 2

# --- Row 15 (new) ---
$ws.Range("A15").Value = "56adfd0fdabb2b70021d45e4"
$ws.Range("B15").Value = "Zoe"
$ws.Range("C15").Value = 56
$ws.Range("D15").Value = "female"
$ws.Range("E15").Value = "no"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 9
$ws.Range("L15").Value = 14502816444676400
$ws.Range("L15").NumberFormat = "#,##0"
$ws.Range("M15").Value = " 296.8372381064871 130.61823889713892 237.86549413094573"
$ws.Range("N15").Value = 2269314070043620
$ws.Range("N15").NumberFormat = "#,##0"
$ws.Range("O15").Value = 17
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = "rfr"
$ws.Range("S15").Value = 3
$ws.Range("T15").Value = "rfrfr"
$ws.Range("U15").Value = "rfrfr"
$ws.Range("V15").Value = 21715208450784000
$ws.Range("V15").NumberFormat = "#,##0"
$ws.Range("W15").Value = 20
$ws.Range("X15").Value = 4
$ws.Range("Y15").Value = 4
$ws.Range("Z15").Value = "oooooo"
$ws.Range("AA15").Value = 4
$ws.Range("AB15").Value = "oooooooooooo"
$ws.Range("AC15").Value = "oo"

# --- Column width adjustments (closest achievable approximations of the
#     author's manual resize / bestFit recompute for the touched columns) ---
$ws.Range("A1").EntireColumn.ColumnWidth = 25.666666666666668
$ws.Range("D1").EntireColumn.ColumnWidth = 6.5
$ws.Range("M1").EntireColumn.ColumnWidth = 54.666666666666664
$ws.Range("N1").EntireColumn.ColumnWidth = 19.333333333333332
$ws.Range("AB1").EntireColumn.ColumnWidth = 14.0
